# Adds a new "2021" column (column R) to the right of the existing
# "2020" column (Q) on the active worksheet, mirroring the format of
# the corresponding Q cell in every row, then updates the selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value = 2021
# Row 4
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 58.14349653559799
# Row 5
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 50.405857641278807
# Row 6
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("R6").Value = 65.995789757646122
# Row 7
$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("R7").Value = 47.339416388110941
# Row 8
$ws.Range("Q8").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = 44.18457369250482
# Row 9
$ws.Range("Q9").Copy()
$ws.Range("R9").PasteSpecial(-4122)
$ws.Range("R9").Value = 50.379263611270765
# Row 10
$ws.Range("Q10").Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("R10").Value = 54.819947539591084
# Row 11
$ws.Range("Q11").Copy()
$ws.Range("R11").PasteSpecial(-4122)
$ws.Range("R11").Value = 47.679920417302263
# Row 12
$ws.Range("Q12").Copy()
$ws.Range("R12").PasteSpecial(-4122)
$ws.Range("R12").Value = 61.861274529713718
# Row 13
$ws.Range("Q13").Copy()
$ws.Range("R13").PasteSpecial(-4122)
$ws.Range("R13").Value = 36.712395096811576
# Row 14
$ws.Range("Q14").Copy()
$ws.Range("R14").PasteSpecial(-4122)
$ws.Range("R14").Value = 26.872053459579295
# Row 15
$ws.Range("Q15").Copy()
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("R15").Value = 46.638444428499682
# Row 16
$ws.Range("Q16").Copy()
$ws.Range("R16").PasteSpecial(-4122)
$ws.Range("R16").Value = 51.155081745820631
# Row 17
$ws.Range("Q17").Copy()
$ws.Range("R17").PasteSpecial(-4122)
$ws.Range("R17").Value = 43.08338023862634
# Row 18
$ws.Range("Q18").Copy()
$ws.Range("R18").PasteSpecial(-4122)
$ws.Range("R18").Value = 58.934228062068456
# Row 19
$ws.Range("Q19").Copy()
$ws.Range("R19").PasteSpecial(-4122)
$ws.Range("R19").Value = 54.51979816984521
# Row 20
$ws.Range("Q20").Copy()
$ws.Range("R20").PasteSpecial(-4122)
$ws.Range("R20").Value = 52.474443936678909
# Row 21
$ws.Range("Q21").Copy()
$ws.Range("R21").PasteSpecial(-4122)
$ws.Range("R21").Value = 56.519551395440942
# Row 22
$ws.Range("Q22").Copy()
$ws.Range("R22").PasteSpecial(-4122)
$ws.Range("R22").Value = 46.970408642555192
# Row 23
$ws.Range("Q23").Copy()
$ws.Range("R23").PasteSpecial(-4122)
$ws.Range("R23").Value = 27.43769048802011
# Row 24
$ws.Range("Q24").Copy()
$ws.Range("R24").PasteSpecial(-4122)
$ws.Range("R24").Value = 66.104415920267911
# Row 25
$ws.Range("Q25").Copy()
$ws.Range("R25").PasteSpecial(-4122)
$ws.Range("R25").Value = 88.246666265390886
# Row 26
$ws.Range("Q26").Copy()
$ws.Range("R26").PasteSpecial(-4122)
$ws.Range("R26").Value = 71.914698721605745
# Row 27
$ws.Range("Q27").Copy()
$ws.Range("R27").PasteSpecial(-4122)
$ws.Range("R27").Value = 105.10059183863845
# Row 28
$ws.Range("Q28").Copy()
$ws.Range("R28").PasteSpecial(-4122)
$ws.Range("R28").Value = 63.980940123966526
# Row 29
$ws.Range("Q29").Copy()
$ws.Range("R29").PasteSpecial(-4122)
$ws.Range("R29").Value = 55.546587096180644
# Row 30
$ws.Range("Q30").Copy()
$ws.Range("R30").PasteSpecial(-4122)
$ws.Range("R30").Value = 73.505198287622903
# Row 31
$ws.Range("Q31").Copy()
$ws.Range("R31").PasteSpecial(-4122)
$ws.Range("R31").Value = 43.916363725083563
# Row 32
$ws.Range("Q32").Copy()
$ws.Range("R32").PasteSpecial(-4122)
$ws.Range("R32").Value = 40.980198843051781
# Row 33
$ws.Range("Q33").Copy()
$ws.Range("R33").PasteSpecial(-4122)
$ws.Range("R33").Value = 47.015458682814909

$excel.CutCopyMode = $false

# Match the author-recorded selection state after the edit.
$ws.Range("T3").Select()
